$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Requi" + "sitos do Aplicativo de" -> one continuous
#    phrase "Requisitos do Aplicativo de". The _GoBack bookmark that used to
#    sit between those two runs is removed from here (it will be re-created
#    later, next to the "-->" text in the "Seguranca" bullet).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) "virtual." -> "virtual;"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("O Sistema deve permitir que os alunos consultem a biblioteca virtual.", $true, $false, $false, $false, $false, $true, 1, $false, "O Sistema deve permitir que os alunos consultem a biblioteca virtual;", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a new list paragraph right after the "virtual;" bullet:
#    "O sistema deve permitir a priorizacao de tarefas." in red.
# ---------------------------------------------------------------------------
$virtualPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*consultem a biblioteca virtual;*") {
        $virtualPara = $p
        break
    }
}
$virtualPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $virtualPara.Next()
$insStart = $newPara.Range.Start

$part1 = $d.Range($insStart, $insStart)
$part1.InsertAfter("O sistema deve pe")
$part1.Font.Color = 255

$part2 = $d.Range($part1.End, $part1.End)
$part2.InsertAfter("rmitir a prioriza" + [char]0x00E7 + [char]0x00E3 + "o de tarefas.")
$part2.Font.Color = 255

# Colour the paragraph mark too, matching the source formatting.
$newPara.Range.Font.Color = 255

# ---------------------------------------------------------------------------
# 4) "contra-ataques; -&gt;" -> "contra-ataques; " + red "-" "-" ">" ,
#    with the _GoBack bookmark re-created right before the red dashes.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("contra-ataques; ->", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = "contra-ataques; "
$dashPos = $r3.End

$dash1 = $d.Range($dashPos, $dashPos)
$dash1.InsertAfter("-")
$dash1.Font.Color = 255

$dash2 = $d.Range($dash1.End, $dash1.End)
$dash2.InsertAfter("-")
$dash2.Font.Color = 255

$dash3 = $d.Range($dash2.End, $dash2.End)
$dash3.InsertAfter(">")
$dash3.Font.Color = 255

$bmRange = $d.Range($dashPos, $dashPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
